# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell whose status was
# "Ready for handoff" is now "In Translation", and the columns holding that
# text (Overview!E:F and the "Status" column on each language sheet) are
# re-sized to fit the new, shorter value.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview!E2:F4 hold the per-language status ("zh-cn" / "de-de" columns).
$overview.Range("E2:F4").Value = "In Translation"

# zh-cn!C2:C4 and de-de!C2:C4 hold the "Status" column for each language table.
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# Resize the affected columns to fit the new, narrower text.
$overview.Range("E:F").ColumnWidth = 12.5
$zhcn.Range("C:C").ColumnWidth = 12.5
$dede.Range("C:C").ColumnWidth = 12.5
